$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.681.10'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.900.41'
$ws.Range('E3').Value = '  +0.68%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  -0.16%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.21'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.31%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -0.17%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5224'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +8.05%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3771'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -0.26%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07239'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -1.25%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.08'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +3.26%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8946'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07625'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.906.90'
$ws.Range('E13').Value = '  +1.15%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.437'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -0.41%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.86'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +1.17%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  -0.20%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008712'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -0.91%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9990'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '27.713.21'
$ws.Range('E19').Value = '  -0.05%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.44'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -0.51%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.125'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '2.139.58'
$ws.Range('E22').Value = '  -1.17%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.81'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -0.01%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.568'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -0.23%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.97'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.25%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.864'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -1.99%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.159'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +2.47%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.26'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.52%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.43'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -1.14%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.831'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -1.19%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08979'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.824'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +4.34%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.236'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +1.48%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7695'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.606'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02073'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +1.66%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.051'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +2.74%  '
$ws.Range('E39').Value = '  +0.04%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5483'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +0.66%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05269'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +0.51%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.632'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -4.44%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '113.13'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +3.42%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.434'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +1.46%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1503'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -0.70%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4777'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +0.11%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.36'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('E48').Value = '  -0.20%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.612'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('E51').Value = '  -0.79%  '
